# Add a new "CDE Version" column (BI) to Sheet1, matching the NINDS
# version-number commit: header label in row 1, and the value 3
# (formatted with two decimal places) in every data row (2-136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column.
$ws.Range("BI1").Value = "CDE Version"

# Data cells: every row from 2 to 136 gets the version number 3,
# displayed with a "0.00" (two-decimal) number format.
$dataRange = $ws.Range("BI2:BI136")
$dataRange.Value = 3
$dataRange.NumberFormat = "0.00"

# Restore the view/selection state recorded in the saved workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 107
$win.ScrollColumn = 36
[void]$ws.Range("BG138").Select()
